# Update gh-pages to output generated at 456a3b4
# Apply the same set of value updates to both the "展览" and "全部类型"
# worksheets (sheet1 and sheet4), which carry duplicate data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Simple numeric column F updates, keyed by row number.
$fUpdates = @{
    4  = 276
    6  = 553
    7  = 57
    8  = 2015
    11 = 4337
    13 = 279
    14 = 98
    15 = 5
    16 = 111
    17 = 24
    19 = 66
    20 = 3126
    22 = 455
    26 = 83
    29 = 53
    32 = 513
    33 = 1717
    34 = 262
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $fUpdates.Keys) {
        $ws.Cells.Item($row, 6).Value = $fUpdates[$row]
    }

    # Row 10: event was cancelled.
    $ws.Range("C10").Value = "赣州·十万伏特-第八届青年文化综合展览会（取消）"
    $ws.Range("F10").Value = 103
    $ws.Range("G10").Value = "不可售"
}
